# Leave Card update — 12/27/2023 upload
# Shifts the PERIOD dates in the OJT leave-card table forward by one
# "cutoff" (30 days) and fills in the EARNED leave credits that have
# accrued for the newly-elapsed cutoffs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Effectivity date at the top of the form (was blank).
$ws.Range("F3").Value = 44743

# PERIOD column (A) — roll every cutoff date forward by 30 days.
$ws.Range("A11").Value = 44773
$ws.Range("A12").Value = 44804
$ws.Range("A13").Value = 44834
$ws.Range("A14").Value = 44865
$ws.Range("A15").Value = 44895
$ws.Range("A16").Value = 44926
$ws.Range("A18").Value = 44957
$ws.Range("A19").Value = 44985
$ws.Range("A20").Value = 45016
$ws.Range("A21").Value = 45046
$ws.Range("A22").Value = 45077
$ws.Range("A23").Value = 45107
$ws.Range("A24").Value = 45138
$ws.Range("A25").Value = 45169
$ws.Range("A26").Value = 45199
$ws.Range("A27").Value = 45230
$ws.Range("A28").Value = 45260
$ws.Range("A29").Value = 45291

# EARNED column (C) — new cutoffs that now have a credited 1.25 entry.
$ws.Range("C20").Value = 1.25
$ws.Range("C21").Value = 1.25
$ws.Range("C22").Value = 1.25
$ws.Range("C23").Value = 1.25
$ws.Range("C24").Value = 1.25
$ws.Range("C25").Value = 1.25
$ws.Range("C26").Value = 1.25

# Move the active cell / selection to reflect where the user left off.
$ws.Range("B15").Select() | Out-Null
